# Update input templates: split the Haber-Bosch ammonia synthesis unit's
# direct "nh3" output into an intermediate "nh3_raw" stream, and add a new
# "nh3_synthesis_power" unit that converts nh3_raw -> nh3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")
$ws.Activate()

# New row 7: nh3_synthesis_power converts nh3_raw (Input1) into nh3 (Output1)
$ws.Range("A7").Value = "nh3_synthesis_power"
$ws.Range("B7").Value = "object_type"
$ws.Range("C7").Value = "nh3_raw"
$ws.Range("E7").Value = "nh3"

# Row 6 (nh3_synthesizer / Haber_Bosch_reactor): Output1 becomes nh3_raw
$ws.Range("E6").Value = "nh3_raw"

# Match the resolution_output / demand column formatting used elsewhere in
# the table (right aligned).
$ws.Range("AG7").HorizontalAlignment = -4152
$ws.Range("AH7").HorizontalAlignment = -4152

# Extend the structured table and related ranges to include the new row.
$table = $ws.ListObjects.Item("Table16")
$table.Resize($ws.Range("A1:AH7"))

# Extend the data-validation list range for the resolution_output column.
$ws.Range("AG2:AG7").Validation.Delete()
$ws.Range("AG2:AG7").Validation.Add(3, 1, 1, """h, D, W, M, Q, Y""")
$ws.Range("AG2:AG7").Validation.IgnoreBlank = $true
$ws.Range("AG2:AG7").Validation.InCellDropdown = $true
$ws.Range("AG2:AG7").Validation.ShowInput = $true
$ws.Range("AG2:AG7").Validation.ShowError = $true

# Widen column A to fit the new, longer unit name (Excel stores column widths
# rounded to the nearest 1/6 character plus its standard padding, so feed it
# the "ColumnWidth" that round-trips to a stored width of exactly 20).
$ws.Columns.Item(1).ColumnWidth = 19.16666666666667

$ws.Range("F21").Select()
